$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "insurance"
$ws.Range("E2").Value = "Health Insurance"
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 5

# Update the visible selection / scroll position to match the saved view state
$ws.Activate()
$ws.Range("L2").Select()
